{"js": "// Replace the date line and each \"a\u00f7b=\" problem in the practice sheet with\n// the new values from the target revision. Every \"from\" string is unique in\n// the document, so a literal, case-sensitive search-and-replace for each\n// pair is unambiguous.\nconst replacements = [\n  [\"2024-03-13 Wednesday\", \"2024-03-14 Thursday\"],\n  [\"15\u00f79=\", \"31\u00f77=\"],\n  [\"72\u00f74=\", \"48\u00f73=\"],\n  [\"74\u00f79=\", \"64\u00f79=\"],\n  [\"83\u00f78=\", \"18\u00f76=\"],\n  [\"77\u00f75=\", \"60\u00f73=\"],\n  [\"40\u00f76=\", \"44\u00f75=\"],\n  [\"21\u00f78=\", \"19\u00f78=\"],\n  [\"68\u00f76=\", \"75\u00f79=\"],\n  [\"46\u00f79=\", \"39\u00f74=\"],\n  [\"91\u00f77=\", \"95\u00f77=\"],\n  [\"13\u00f79=\", \"34\u00f79=\"],\n  [\"99\u00f72=\", \"73\u00f75=\"],\n  [\"67\u00f79=\", \"25\u00f78=\"],\n  [\"87\u00f72=\", \"62\u00f72=\"],\n  [\"83\u00f76=\", \"17\u00f74=\"],\n  [\"21\u00f76=\", \"95\u00f72=\"],\n  [\"82\u00f75=\", \"64\u00f75=\"],\n  [\"64\u00f77=\", \"79\u00f74=\"],\n  [\"25\u00f77=\", \"41\u00f76=\"],\n  [\"28\u00f72=\", \"70\u00f73=\"],\n  [\"74\u00f74=\", \"74\u00f78=\"],\n  [\"18\u00f72=\", \"50\u00f77=\"],\n  [\"87\u00f77=\", \"25\u00f78=\"],\n  [\"36\u00f73=\", \"35\u00f75=\"],\n  [\"52\u00f78=\", \"84\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const results = body.search(from, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"a\u00f7b=\" problem in the practice sheet with\n# the new values from the target revision. Every \"from\" string is unique in\n# the document, so a Find/Replace scoped to the whole story for each pair is\n# unambiguous.\n$replacements = @(\n    @('2024-03-13 Wednesday', '2024-03-14 Thursday'),\n    @('15\u00f79=', '31\u00f77='),\n    @('72\u00f74=', '48\u00f73='),\n    @('74\u00f79=', '64\u00f79='),\n    @('83\u00f78=', '18\u00f76='),\n    @('77\u00f75=', '60\u00f73='),\n    @('40\u00f76=', '44\u00f75='),\n    @('21\u00f78=', '19\u00f78='),\n    @('68\u00f76=', '75\u00f79='),\n    @('46\u00f79=', '39\u00f74='),\n    @('91\u00f77=', '95\u00f77='),\n    @('13\u00f79=', '34\u00f79='),\n    @('99\u00f72=', '73\u00f75='),\n    @('67\u00f79=', '25\u00f78='),\n    @('87\u00f72=', '62\u00f72='),\n    @('83\u00f76=', '17\u00f74='),\n    @('21\u00f76=', '95\u00f72='),\n    @('82\u00f75=', '64\u00f75='),\n    @('64\u00f77=', '79\u00f74='),\n    @('25\u00f77=', '41\u00f76='),\n    @('28\u00f72=', '70\u00f73='),\n    @('74\u00f74=', '74\u00f78='),\n    @('18\u00f72=', '50\u00f77='),\n    @('87\u00f77=', '25\u00f78='),\n    @('36\u00f73=', '35\u00f75='),\n    @('52\u00f78=', '84\u00f73=')\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
